# Add a "Save" column (H) to the s_vals sheet, mirroring the header style
# used by the existing "sum" column (G) and filling the data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same bold/bordered/centered style as the other headers (G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells: new "Save" values, defaulting to 0.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
